$d = $word.ActiveDocument

$rFonts = '<w:rFonts w:ascii="Book Antiqua" w:eastAsia="Book Antiqua" w:hAnsi="Book Antiqua" w:cs="Book Antiqua"/>'
$color  = '<w:color w:val="000000"/>'
$pPrCommon = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:hanging="851"/><w:rPr>' + $rFonts + $color + '</w:rPr></w:pPr>'

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Locate the "c." and "d." declaration paragraphs by their current wording.
$cPara = $null
$dPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*No export incentive has been availed*") { $cPara = $p }
    if ($t -like "*Proportionate export incentives shall be surrendered*") { $dPara = $p }
}

# --- "c." paragraph: re-run the grammar check so "case  OR" becomes its own run. ---
# Only the runs (not the paragraph mark) are replaced, so the existing <w:p>
# attributes (w14:paraId, rsids, ...) and <w:pPr> survive untouched.
$cFull = $cPara.Range
$cInner = $d.Range($cFull.Start, $cFull.End - 1)
$cRuns = ''
$cRuns += '<w:r><w:rPr>' + $rFonts + '<w:b/>' + $color + '</w:rPr><w:t>c.</w:t></w:r>'
$cRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t xml:space="preserve">    No export incentive has been availed for the reductions / deductions. If that is the </w:t></w:r>'
$cRuns += '<w:proofErr w:type="gramStart"/>'
$cRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t>case  OR</w:t></w:r>'
$cRuns += '<w:proofErr w:type="gramEnd"/>'
$cRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t xml:space="preserve">  Proportionate export incentives have been surrendered If that is the case  OR</w:t></w:r>'
$cInner.InsertXML((New-PkgXml ('<w:p>' + $cRuns + '</w:p>')))

# --- "d." paragraph: replace the wording with the new declaration text. ---
$dFull = $dPara.Range
$dInner = $d.Range($dFull.Start, $dFull.End - 1)
$dRuns = ''
$dRuns += '<w:r><w:rPr>' + $rFonts + '<w:b/>' + $color + '</w:rPr><w:t>d.</w:t></w:r>'
$dRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
$dRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t xml:space="preserve">  </w:t></w:r>'
$dRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t>No proportionate export incentives availed</w:t></w:r>'
$dRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t xml:space="preserve">. </w:t></w:r>'
$dInner.InsertXML((New-PkgXml ('<w:p>' + $dRuns + '</w:p>')))

# --- Insert the new "e." declaration paragraph right after "d." ---
$dPara.Range.InsertParagraphAfter()
$ePara = $dPara.Next()
$eRuns = ''
$eRuns += '<w:r><w:rPr>' + $rFonts + '<w:b/>' + $color + '</w:rPr><w:t>e.</w:t></w:r>'
$eRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t xml:space="preserve">   </w:t></w:r>'
$eRuns += '<w:proofErr w:type="spellStart"/>'
$eRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t>ChicMic</w:t></w:r>'
$eRuns += '<w:proofErr w:type="spellEnd"/>'
$eRuns += '<w:r><w:rPr>' + $rFonts + $color + '</w:rPr><w:t xml:space="preserve"> Technologies LLP is not in Caution List</w:t></w:r>'
# Use the full paragraph range (it is brand new/empty) so the fresh <w:pPr>
# has to be supplied explicitly, but no stray empty run is left behind.
$ePara.Range.InsertXML((New-PkgXml ('<w:p>' + $pPrCommon + $eRuns + '</w:p>')))

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
